# Applies the "Algorithms tests comparison table" edit:
#  - adds a new corner header label in A1 ("Algorithm / Action") with a
#    diagonal divider border
#  - switches the header/row-label font from Calibri to Cambria (theme major
#    font) and re-centers the top header row (drops the 45-degree rotation)
#  - re-percent-formats + re-centers the data cells in Cambria (regular),
#    adds a new 7th results column (G) and widens columns B:G uniformly
#  - tweaks row heights for the two data rows and moves the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new corner header label (A1) ------------------------------------------
$cornerText = @"
    Algorithm     
                                 Action
"@
$ws.Range("A1").Value = $cornerText

# ---- new 7th data column (G) -----------------------------------------------
$ws.Range("G2").Value = 0.90910000000000002
$ws.Range("G3").Value = 0.91669999999999996

# ---- column widths: B:G become one uniform width ---------------------------
$ws.Range("B1:G1").EntireColumn.ColumnWidth = 12.43

# ---- row heights for the two data rows -------------------------------------
$ws.Rows.Item(2).RowHeight = 43.5
$ws.Rows.Item(3).RowHeight = 42.75

# ---- A1: Cambria bold, diagonal divider, left/center/wrap ------------------
$a1 = $ws.Range("A1")
$a1.Font.Name = "Cambria"
$a1.Font.Bold = $true
$a1.Borders.Item(6).LineStyle = 1
$a1.Borders.Item(6).Weight = 2
$a1.HorizontalAlignment = -4131
$a1.VerticalAlignment = -4108
$a1.WrapText = $true

# ---- B1:G1 header cells: Cambria bold, centered (no more 45-deg rotation) --
$headerRow = $ws.Range("B1:G1")
$headerRow.Font.Name = "Cambria"
$headerRow.Font.Bold = $true
$headerRow.Orientation = 0
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4108
$headerRow.WrapText = $true

# ---- A2:A3 row labels: Cambria bold, left/center ---------------------------
$rowLabels = $ws.Range("A2:A3")
$rowLabels.Font.Name = "Cambria"
$rowLabels.Font.Bold = $true
$rowLabels.HorizontalAlignment = -4131
$rowLabels.VerticalAlignment = -4108
$rowLabels.WrapText = $true

# ---- B2:G3 data cells: Cambria regular, centered, percent formats ----------
$dataRange = $ws.Range("B2:G3")
$dataRange.Font.Name = "Cambria"
$dataRange.Font.Bold = $false
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true
$dataRange.NumberFormat = "0.00%"

# C2:D2 keep the "0%" (no-decimal) percent format
$ws.Range("C2:D2").NumberFormat = "0%"

# ---- move the active selection to A2:A3, as in the edited file ------------
$ws.Range("A2:A3").Select() | Out-Null
